$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair results (Natmi, following Dr Hou advice): adds an "ECs"
# sending-cluster and expands the row set to cover every Sending x Target
# cluster combination (ECs/FAPs/sCs) for the Vtn -> Tnfrsf11b pair.

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vtn"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "FAPs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 7.134618
$ws.Cells.Item(2, 8).Value = 21.403854
$ws.Cells.Item(2, 9).Value = 0.0965317920926077
$ws.Cells.Item(2, 10).Value = 0.0965317920926077
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.095195666666667
$ws.Cells.Item(2, 14).Value = 6.285587
$ws.Cells.Item(2, 15).Value = 0.8546922300706357
$ws.Cells.Item(2, 16).Value = 0.8546922300706358
$ws.Cells.Item(2, 17).Value = 14.948420716922
$ws.Cells.Item(2, 18).Value = 134.535786452298
$ws.Cells.Item(2, 19).Value = 0.08250497265634583
$ws.Cells.Item(2, 20).Value = 0.08250497265634585

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vtn"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "sCs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 7.134618
$ws.Cells.Item(3, 8).Value = 21.403854
$ws.Cells.Item(3, 9).Value = 0.0965317920926077
$ws.Cells.Item(3, 10).Value = 0.0965317920926077
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.356208
$ws.Cells.Item(3, 14).Value = 1.068624
$ws.Cells.Item(3, 15).Value = 0.1453077699293643
$ws.Cells.Item(3, 16).Value = 0.1453077699293643
$ws.Cells.Item(3, 17).Value = 2.541408008544
$ws.Cells.Item(3, 18).Value = 22.872672076896
$ws.Cells.Item(3, 19).Value = 0.01402681943626187
$ws.Cells.Item(3, 20).Value = 0.01402681943626187

$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Vtn"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "FAPs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 17.50798033333334
$ws.Cells.Item(4, 8).Value = 52.52394100000001
$ws.Cells.Item(4, 9).Value = 0.2368839813846793
$ws.Cells.Item(4, 10).Value = 0.2368839813846794
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.095195666666667
$ws.Cells.Item(4, 14).Value = 6.285587
$ws.Cells.Item(4, 15).Value = 0.8546922300706357
$ws.Cells.Item(4, 16).Value = 0.8546922300706358
$ws.Cells.Item(4, 17).Value = 36.68264452648523
$ws.Cells.Item(4, 18).Value = 330.143800738367
$ws.Cells.Item(4, 19).Value = 0.2024628983176825
$ws.Cells.Item(4, 20).Value = 0.2024628983176826

$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Vtn"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 17.50798033333334
$ws.Cells.Item(5, 8).Value = 52.52394100000001
$ws.Cells.Item(5, 9).Value = 0.2368839813846793
$ws.Cells.Item(5, 10).Value = 0.2368839813846794
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.356208
$ws.Cells.Item(5, 14).Value = 1.068624
$ws.Cells.Item(5, 15).Value = 0.1453077699293643
$ws.Cells.Item(5, 16).Value = 0.1453077699293643
$ws.Cells.Item(5, 17).Value = 6.236482658576001
$ws.Cells.Item(5, 18).Value = 56.12834392718401
$ws.Cells.Item(5, 19).Value = 0.0344210830669968
$ws.Cells.Item(5, 20).Value = 0.0344210830669968

$ws.Cells.Item(6, 1).Value = "sCs"
$ws.Cells.Item(6, 2).Value = "Vtn"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 49.26691733333334
$ws.Cells.Item(6, 8).Value = 147.800752
$ws.Cells.Item(6, 9).Value = 0.6665842265227129
$ws.Cells.Item(6, 10).Value = 0.666584226522713
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.095195666666667
$ws.Cells.Item(6, 14).Value = 6.285587
$ws.Cells.Item(6, 15).Value = 0.8546922300706357
$ws.Cells.Item(6, 16).Value = 0.8546922300706358
$ws.Cells.Item(6, 17).Value = 103.2238317068249
$ws.Cells.Item(6, 18).Value = 929.014485361424
$ws.Cells.Item(6, 19).Value = 0.5697243590966072
$ws.Cells.Item(6, 20).Value = 0.5697243590966073

$ws.Cells.Item(7, 1).Value = "sCs"
$ws.Cells.Item(7, 2).Value = "Vtn"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 49.26691733333334
$ws.Cells.Item(7, 8).Value = 147.800752
$ws.Cells.Item(7, 9).Value = 0.6665842265227129
$ws.Cells.Item(7, 10).Value = 0.666584226522713
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.356208
$ws.Cells.Item(7, 14).Value = 1.068624
$ws.Cells.Item(7, 15).Value = 0.1453077699293643
$ws.Cells.Item(7, 16).Value = 0.1453077699293643
$ws.Cells.Item(7, 17).Value = 17.549270089472
$ws.Cells.Item(7, 18).Value = 157.943430805248
$ws.Cells.Item(7, 19).Value = 0.09685986742610561
$ws.Cells.Item(7, 20).Value = 0.09685986742610563
